$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the distribution values for row 7 (9 players):
#  C7: 3 -> 2
#  D7: "-" -> 1
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1

# Move/leave the active selection on D13, matching the author's last
# selected cell when the workbook was saved.
[void]$ws.Range("D13").Select()
